# "fixed harvester column in rnasamples -- holly added S.GISH to harvester in bioSamples"
#
# The 'harvester' column (column B) should read "S.GISH" for every data
# row instead of the old "Retrofitted_0659" placeholder value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the harvester column, as the editor would have done before typing
# the replacement value in.
$ws.Range("B:B").Select()

# Data rows are 2 through 21 (row 1 is the header row).
$ws.Range("B2:B21").Value = "S.GISH"
